$wb = $excel.ActiveWorkbook

# Map of worksheet name -> new F2/F3 values (Excel date serial numbers)
$updates = @{
    "assets"       = @{ F2 = 43104; F3 = 43101 }
    "Assets2"      = @{ F2 = 43104; F3 = 43101 }
    "Liabilities"  = @{ F2 = 43101; F3 = 43102 }
    "Liabilities2" = @{ F2 = 43101; F3 = 43102 }
    "Assets3"      = @{ F2 = 43104; F3 = 43101 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $vals = $updates[$sheetName]
    $ws.Range("F2").Value = $vals.F2
    $ws.Range("F3").Value = $vals.F3
}
